# Auto-generated edit script applying market-price refresh values
# as described by the commit diff, per sheet/cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4107.2856
$ws.Range("I62").Value = 3291.8333
$ws.Range("K62").Value = 3291.8333
$ws.Range("M62").Value = -2667.8333
$ws.Range("H65").Value = 4107.2856
$ws.Range("I65").Value = 3291.8333
$ws.Range("K65").Value = 16459.1665
$ws.Range("M65").Value = -13339.1665
$ws.Range("H70").Value = 2256.121
$ws.Range("I70").Value = 1418.8889
$ws.Range("K70").Value = 4256.6667
$ws.Range("M70").Value = -3986.6667
$ws.Range("H73").Value = 2256.121
$ws.Range("I73").Value = 1418.8889
$ws.Range("K73").Value = 4256.6667
$ws.Range("M73").Value = -3320.6667
$ws.Range("H74").Value = 10944.333
$ws.Range("I74").Value = 9499.666999999999
$ws.Range("K74").Value = 9499.666999999999
$ws.Range("M74").Value = -8563.666999999999
$ws.Range("H77").Value = 10944.333
$ws.Range("I77").Value = 9499.666999999999
$ws.Range("K77").Value = 47498.335
$ws.Range("M77").Value = -42818.335
$ws.Range("H98").Value = 1039.3793
$ws.Range("I98").Value = 1141
$ws.Range("J98").Value = 720
$ws.Range("K98").Value = 1141
$ws.Range("L98").Value = 720
$ws.Range("M98").Value = 357
$ws.Range("N98").Value = -3716
$ws.Range("H112").Value = 1885.6897
$ws.Range("I112").Value = 1828.3334
$ws.Range("J112").Value = 1892.3077
$ws.Range("K112").Value = 5485.0002
$ws.Range("L112").Value = 5676.9231
$ws.Range("M112").Value = -4377.0002
$ws.Range("N112").Value = -7892.9231
$ws.Range("H116").Value = 4439.2
$ws.Range("I116").Value = 4624.25
$ws.Range("J116").Value = 3699
$ws.Range("K116").Value = 4624.25
$ws.Range("L116").Value = 3699
$ws.Range("M116").Value = -1182.25
$ws.Range("N116").Value = -10583
$ws.Range("H122").Value = 1039.3793
$ws.Range("I122").Value = 1141
$ws.Range("J122").Value = 720
$ws.Range("K122").Value = 3423
$ws.Range("L122").Value = 2160
$ws.Range("M122").Value = -973
$ws.Range("N122").Value = -7060
$ws.Range("H129").Value = 126827.875
$ws.Range("I129").Value = 183395.36
$ws.Range("J129").Value = 2379.4
$ws.Range("K129").Value = 550186.08
$ws.Range("L129").Value = 7138.200000000001
$ws.Range("M129").Value = -545186.08
$ws.Range("N129").Value = -17138.2
$ws.Range("H130").Value = 65000
$ws.Range("J130").Value = 65000
$ws.Range("L130").Value = 65000
$ws.Range("N130").Value = -75040
$ws.Range("H137").Value = 2412.4036
$ws.Range("I137").Value = 2273.9788
$ws.Range("K137").Value = 6821.9364
$ws.Range("M137").Value = -4271.9364
$ws.Range("H138").Value = 3022.1614
$ws.Range("I138").Value = 1278.5333
$ws.Range("J138").Value = 4656.8125
$ws.Range("K138").Value = 3835.5999
$ws.Range("L138").Value = 13970.4375
$ws.Range("M138").Value = 1304.4001
$ws.Range("N138").Value = -24250.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16957332
$ws.Range("I32").Value = 17865440
$ws.Range("J32").Value = 5999.3335
$ws.Range("K32").Value = 17865440
$ws.Range("L32").Value = 5999.3335
$ws.Range("M32").Value = -17865153
$ws.Range("N32").Value = -6573.3335
$ws.Range("H74").Value = 2459.2903
$ws.Range("J74").Value = 2100.111
$ws.Range("L74").Value = 2100.111
$ws.Range("N74").Value = -3848.111
$ws.Range("H77").Value = 2459.2903
$ws.Range("J77").Value = 2100.111
$ws.Range("L77").Value = 10500.555
$ws.Range("N77").Value = -19236.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H64").Value = 972.44446
$ws.Range("I64").Value = 818.3333
$ws.Range("J64").Value = 1049.5
$ws.Range("K64").Value = 818.3333
$ws.Range("L64").Value = 1049.5
$ws.Range("M64").Value = -593.3333
$ws.Range("N64").Value = -1499.5
$ws.Range("H67").Value = 972.44446
$ws.Range("I67").Value = 818.3333
$ws.Range("J67").Value = 1049.5
$ws.Range("K67").Value = 818.3333
$ws.Range("L67").Value = 1049.5
$ws.Range("M67").Value = -38.33330000000001
$ws.Range("N67").Value = -2609.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2381.8333
$ws.Range("J31").Value = 3484.5
$ws.Range("L31").Value = 3484.5
$ws.Range("N31").Value = -4074.5
$ws.Range("H34").Value = 2381.8333
$ws.Range("J34").Value = 3484.5
$ws.Range("L34").Value = 3484.5
$ws.Range("N34").Value = -3888.5
$ws.Range("H107").Value = 13288.4375
$ws.Range("I107").Value = 635.9091
$ws.Range("J107").Value = 41124
$ws.Range("K107").Value = 635.9091
$ws.Range("L107").Value = 41124
$ws.Range("M107").Value = 1284.0909
$ws.Range("N107").Value = -44964
$ws.Range("H134").Value = 5000.933
$ws.Range("I134").Value = 4134.5557
$ws.Range("K134").Value = 12403.6671
$ws.Range("M134").Value = -9868.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1249.9445
$ws.Range("I113").Value = 1405.75
$ws.Range("J113").Value = 1205.4286
$ws.Range("K113").Value = 4217.25
$ws.Range("L113").Value = 3616.2858
$ws.Range("M113").Value = -2047.25
$ws.Range("N113").Value = -7956.2858
$ws.Range("H129").Value = 1333.1428
$ws.Range("J129").Value = 2114.8333
$ws.Range("L129").Value = 6344.499899999999
$ws.Range("N129").Value = -16344.4999
$ws.Range("H134").Value = 9359.23
$ws.Range("I134").Value = 3381.4285
$ws.Range("K134").Value = 10144.2855
$ws.Range("M134").Value = -5074.2855
$ws.Range("H138").Value = 7218.4
$ws.Range("I138").Value = 6498.5
$ws.Range("K138").Value = 19495.5
$ws.Range("M138").Value = -14355.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3662.5862
$ws.Range("I122").Value = 1262.6
$ws.Range("J122").Value = 4925.737
$ws.Range("K122").Value = 3787.8
$ws.Range("L122").Value = 14777.211
$ws.Range("M122").Value = -1337.8
$ws.Range("N122").Value = -19677.211

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1899.8182
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 2066.3333
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 2066.3333
$ws.Range("M22").Value = -1405
$ws.Range("N22").Value = -2656.3333
$ws.Range("H27").Value = 1899.8182
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 2066.3333
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 2066.3333
$ws.Range("M27").Value = -1593
$ws.Range("N27").Value = -2280.3333
$ws.Range("H55").Value = 6857.143
$ws.Range("I55").Value = 7433
$ws.Range("K55").Value = 7433
$ws.Range("M55").Value = -7260
$ws.Range("H59").Value = 6649.6665
$ws.Range("I59").Value = 5000
$ws.Range("J59").Value = 7474.5
$ws.Range("K59").Value = 5000
$ws.Range("L59").Value = 7474.5
$ws.Range("M59").Value = -4346
$ws.Range("N59").Value = -8782.5
$ws.Range("H61").Value = 3515.3333
$ws.Range("J61").Value = 4994.5
$ws.Range("L61").Value = 4994.5
$ws.Range("N61").Value = -5398.5
$ws.Range("H68").Value = 3941.3333
$ws.Range("I68").Value = 3943.25
$ws.Range("J68").Value = 3939.8
$ws.Range("K68").Value = 3943.25
$ws.Range("L68").Value = 3939.8
$ws.Range("M68").Value = -3194.25
$ws.Range("N68").Value = -5437.8
$ws.Range("H69").Value = 76719.336
$ws.Range("J69").Value = 76719.336
$ws.Range("L69").Value = 76719.336
$ws.Range("N69").Value = -78341.336
$ws.Range("H71").Value = 3941.3333
$ws.Range("I71").Value = 3943.25
$ws.Range("J71").Value = 3939.8
$ws.Range("K71").Value = 19716.25
$ws.Range("L71").Value = 19699
$ws.Range("M71").Value = -15972.25
$ws.Range("N71").Value = -27187
$ws.Range("H72").Value = 76719.336
$ws.Range("J72").Value = 76719.336
$ws.Range("L72").Value = 230158.008
$ws.Range("N72").Value = -238270.008
$ws.Range("H113").Value = 3515.3333
$ws.Range("J113").Value = 4994.5
$ws.Range("L113").Value = 4994.5
$ws.Range("N113").Value = -9334.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 19964.8
$ws.Range("J45").Value = 22306.25
$ws.Range("L45").Value = 22306.25
$ws.Range("N45").Value = -23288.25
$ws.Range("H113").Value = 859.8823
$ws.Range("I113").Value = 558.6429000000001
$ws.Range("J113").Value = 2265.6667
$ws.Range("K113").Value = 1675.9287
$ws.Range("L113").Value = 6797.000100000001
$ws.Range("M113").Value = 494.0712999999998
$ws.Range("N113").Value = -11137.0001
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H124").Value = 234999.75
$ws.Range("J124").Value = 234999.75
$ws.Range("L124").Value = 234999.75
$ws.Range("N124").Value = -244819.75
$ws.Range("H126").Value = 1085.5264
$ws.Range("I126").Value = 1090.2778
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 3270.8334
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -800.8334000000004
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 1868.3
$ws.Range("I132").Value = 1830.7567
$ws.Range("K132").Value = 5492.2701
$ws.Range("M132").Value = -2962.2701

Write-Host "Applied 235 cell updates and 2 cell clears."